# Append new order line items (rows 23-38) to the active worksheet.
# Source data mirrors the existing sheet layout:
#   A = SKU, B = Name, C = Quantity, D = Cost Per, E = Total Cost
# All values in the existing sheet are stored as text, so each new value
# is written with a leading apostrophe (forces text) and then the cell's
# style is reset back to Normal so no numeric formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$rows = @(
    @("6422273", "Lemon Juice",               "1",  "16.55", "16.55"),
    @("8255796", "Tuna White Chunk (Pouch)",  "6",  "72.00", "432.00"),
    @("1132582", "Sour Cream",                "1",  "28.94", "28.94"),
    @("7722184", "Parmesan (Grated)",         "1",  "59.95", "59.95"),
    @("1027629", "Cheddar - (Sliced)",        "12", "34.91", "418.92"),
    @("1035842", "Feta - Pail",               "1",  "92.87", "92.87"),
    @("6364494", "Yogurt - Greek (Bulk)",     "2",  "27.11", "54.22"),
    @("4254710", "Veggie Burger",             "1",  "52.72", "52.72"),
    @("1365278", "Vegan Chicken Tenders",     "1",  "87.80", "87.80"),
    @("3275539", "Sauerkraut",                "3",  "19.35", "58.05"),
    @("1028188", "Tortellini - Cheese",       "2",  "32.22", "64.44"),
    @("7529232", 'Wrap - Wheat (10")',         "2",  "31.32", "62.64"),
    @("2825368", "Sausage - Chicken Patty",   "3",  "50.01", "150.03"),
    @("4157160", "Spanakopita",               "2",  "75.51", "151.02"),
    @("9546982", "Arugula - Fresh",           "4",  "20.35", "81.40"),
    @("6264014", "Tomato - Grape",            "6",  "18.06", "108.36")
)

$startRow = 23
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $colNum = $c + 1
        $val = $data[$c]
        Set-TextCell $r $colNum $val
    }
}
